# Recibo template: drop the explicit 12pt (sz/szCs=24) direct formatting from
# the "TOTAL" amount cell ($${Recibo.Total}.00) while keeping it bold, per the
# commit "Formato del plantilla recibos".
#
# The Word COM layer here always serialises Font.Size as a literal half-point
# value (there's no "unset"/inherit sentinel), so the only faithful way to
# drop the w:sz/w:szCs attributes entirely is to replace the paragraph's raw
# OOXML in place via Range.InsertXML with the same markup minus those two
# elements.

$d = $word.ActiveDocument

$targetPara = $null
foreach ($t in $d.Tables) {
    foreach ($row in $t.Rows) {
        foreach ($cell in $row.Cells) {
            if ($cell.Range.Text -like "*Recibo.Total*") {
                $targetPara = $cell.Range.Paragraphs.Item(1)
            }
        }
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate the Recibo.Total paragraph"
}

$paraRange = $targetPara.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' + `
       'w14:paraId="1991B09A" w14:textId="5EBC2D3D" w:rsidR="007B0E30" ' + `
       'w:rsidRPr="005013F0" w:rsidRDefault="007B0E30" w:rsidP="007B0E30">' + `
         '<w:pPr>' + `
           '<w:jc w:val="center"/>' + `
           '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
         '</w:pPr>' + `
         '<w:r w:rsidRPr="005013F0">' + `
           '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
           '<w:t>$${</w:t>' + `
         '</w:r>' + `
         '<w:proofErr w:type="spellStart"/>' + `
         '<w:r w:rsidRPr="005013F0">' + `
           '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
           '<w:t>Recibo.Total</w:t>' + `
         '</w:r>' + `
         '<w:proofErr w:type="spellEnd"/>' + `
         '<w:r w:rsidRPr="005013F0">' + `
           '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
           '<w:t>}</w:t>' + `
         '</w:r>' + `
         '<w:r w:rsidR="00D1180E" w:rsidRPr="005013F0">' + `
           '<w:rPr><w:b/><w:bCs/></w:rPr>' + `
           '<w:t>.00</w:t>' + `
         '</w:r>' + `
       '</w:p>'

$paraRange.InsertXML($xml)
